$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows right above the old row 699 (current weekly snapshot
# for "Betarraga" gets added at the top of the historical table; every row
# below shifts down by two, dimension grows from R734 to R736).
$ws.Rows.Item(699).Insert()
$ws.Rows.Item(699).Insert()

# New row 699: "Primera" quality entry for the new date.
$ws.Cells.Item(699,1).Value = 3
$ws.Cells.Item(699,2).Value = "Femacal de La Calera"
$ws.Cells.Item(699,3).Value = "Coquimbo"
$ws.Cells.Item(699,4).Value = 44753
$ws.Cells.Item(699,5).Value = 5
$ws.Cells.Item(699,6).Value = 100114014
$ws.Cells.Item(699,7).Value = "Betarraga"
$ws.Cells.Item(699,8).Value = "Sin especificar"
$ws.Cells.Item(699,9).Value = "Primera"
$ws.Cells.Item(699,10).Value = 3700
$ws.Cells.Item(699,11).Value = 800
$ws.Cells.Item(699,12).Value = 900
$ws.Cells.Item(699,13).Value = 849
$ws.Cells.Item(699,14).Value = "$/paquete 4 unidades"
$ws.Cells.Item(699,15).Value = "Provincia de Quillota"
$ws.Cells.Item(699,16).Value = 212
$ws.Cells.Item(699,17).Value = 4
$ws.Cells.Item(699,18).Value = "Hortaliza"

# New row 700: "Segunda" quality entry for the same new date.
$ws.Cells.Item(700,1).Value = 3
$ws.Cells.Item(700,2).Value = "Femacal de La Calera"
$ws.Cells.Item(700,3).Value = "Coquimbo"
$ws.Cells.Item(700,4).Value = 44753
$ws.Cells.Item(700,5).Value = 5
$ws.Cells.Item(700,6).Value = 100114014
$ws.Cells.Item(700,7).Value = "Betarraga"
$ws.Cells.Item(700,8).Value = "Sin especificar"
$ws.Cells.Item(700,9).Value = "Segunda"
$ws.Cells.Item(700,10).Value = 1300
$ws.Cells.Item(700,11).Value = 600
$ws.Cells.Item(700,12).Value = 600
$ws.Cells.Item(700,13).Value = 600
$ws.Cells.Item(700,14).Value = "$/paquete 4 unidades"
$ws.Cells.Item(700,15).Value = "Provincia de Quillota"
$ws.Cells.Item(700,16).Value = 150
$ws.Cells.Item(700,17).Value = 4
$ws.Cells.Item(700,18).Value = "Hortaliza"
